$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last data row (96) into the new row (97),
# matching the style pattern used throughout the sheet (index column bold
# bordered style, date column date-time number format, rest General).
$ws.Range("A96:V96").Copy()
$ws.Range("A97:V97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = 97

$ws.Cells.Item($row, 1).Value = 96
$ws.Cells.Item($row, 2).Value = "netherlands"
$ws.Cells.Item($row, 3).Value = "tweede-divisie"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45242.58333333334
$ws.Cells.Item($row, 6).Value = "ADO 20 Heemskerk"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Lisse"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.49
$ws.Cells.Item($row, 11).Value = "12/11/2023 11:12"
$ws.Cells.Item($row, 12).Value = 1.54
$ws.Cells.Item($row, 13).Value = "12/11/2023 13:46"
$ws.Cells.Item($row, 14).Value = 4.77
$ws.Cells.Item($row, 15).Value = "12/11/2023 11:12"
$ws.Cells.Item($row, 16).Value = 4.57
$ws.Cells.Item($row, 17).Value = "12/11/2023 13:48"
$ws.Cells.Item($row, 18).Value = 4.55
$ws.Cells.Item($row, 19).Value = "12/11/2023 11:12"
$ws.Cells.Item($row, 20).Value = 4.63
$ws.Cells.Item($row, 21).Value = "12/11/2023 13:47"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/netherlands/tweede-divisie/ado-20-heemskerk-lisse/jo7NCy0g/"
